$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 680
$ws.Range("E2").Value = 335
$ws.Range("F2").Value = 335
$ws.Range("G2").Value = 313
$ws.Range("H2").Value = 238
$ws.Range("I2").Value = 238
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 2854
$ws.Range("L2").Value = 1048
$ws.Range("M2").Value = 1806
$ws.Range("N2").Value = 1806
$ws.Range("P2").Value = 355
$ws.Range("Q2").Value = 27
$ws.Range("R2").Value = -174
$ws.Range("S2").Value = 219
$ws.Range("T2").Value = 5
$ws.Range("V2").Value = 509
$ws.Range("W2").Value = 49.32
$ws.Range("X2").Value = 35.01
$ws.Range("Y2").Value = 14.83
$ws.Range("Z2").Value = 9.369999999999999
$ws.Range("AA2").Value = 58.02
$ws.Range("AB2").Value = 408.69
$ws.Range("AC2").Value = 302
$ws.Range("AE2").Value = 1744
$ws.Range("AF2").Value = 0
$ws.Range("AG2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 103572356
$ws.Range("O2").ClearContents()
$ws.Range("U2").ClearContents()
$ws.Range("AD2").ClearContents()
$ws.Range("AH2").ClearContents()

# Row 3
$ws.Range("D3").Value = 954
$ws.Range("E3").Value = 584
$ws.Range("F3").Value = 584
$ws.Range("G3").Value = 569
$ws.Range("H3").Value = 422
$ws.Range("I3").Value = 423
$ws.Range("J3").Value = -1
$ws.Range("K3").Value = 3507
$ws.Range("L3").Value = 1170
$ws.Range("M3").Value = 2337
$ws.Range("N3").Value = 2308
$ws.Range("O3").Value = 29
$ws.Range("P3").Value = 367
$ws.Range("Q3").Value = -2
$ws.Range("R3").Value = 236
$ws.Range("S3").Value = -63
$ws.Range("T3").Value = 5
$ws.Range("V3").Value = 388
$ws.Range("W3").Value = 61.2
$ws.Range("X3").Value = 44.23
$ws.Range("Y3").Value = 20.57
$ws.Range("Z3").Value = 13.27
$ws.Range("AA3").Value = 50.05
$ws.Range("AB3").Value = 536.05
$ws.Range("AC3").Value = 406
$ws.Range("AE3").Value = 2153
$ws.Range("AF3").Value = 0
$ws.Range("AG3").Value = 120
$ws.Range("AI3").Value = 22.07
$ws.Range("AJ3").Value = 107186089
$ws.Range("U3").ClearContents()
$ws.Range("AD3").ClearContents()
$ws.Range("AH3").ClearContents()

# Row 4
$ws.Range("D4").Value = 1365
$ws.Range("E4").Value = 962
$ws.Range("F4").Value = 962
$ws.Range("G4").Value = 960
$ws.Range("H4").Value = 732
$ws.Range("I4").Value = 728
$ws.Range("J4").Value = 4
$ws.Range("K4").Value = 5676
$ws.Range("L4").Value = 1506
$ws.Range("M4").Value = 4170
$ws.Range("N4").Value = 4137
$ws.Range("O4").Value = 32
$ws.Range("P4").Value = 425
$ws.Range("Q4").Value = -1043
$ws.Range("R4").Value = -425
$ws.Range("S4").Value = 1218
$ws.Range("T4").Value = 2
$ws.Range("V4").Value = 511
$ws.Range("W4").Value = 70.48
$ws.Range("X4").Value = 53.61
$ws.Range("Y4").Value = 22.59
$ws.Range("Z4").Value = 15.94
$ws.Range("AA4").Value = 36.12
$ws.Range("AB4").Value = 881.12
$ws.Range("AC4").Value = 631
$ws.Range("AD4").Value = 9.94
$ws.Range("AE4").Value = 3337
$ws.Range("AF4").Value = 2.06
$ws.Range("AG4").Value = 137
$ws.Range("AH4").Value = 2.19
$ws.Range("AI4").Value = 23.35
$ws.Range("AJ4").Value = 123977752
$ws.Range("U4").ClearContents()

# Row 5
$ws.Range("D5").Value = 2225
$ws.Range("E5").Value = 1668
$ws.Range("F5").Value = 1668
$ws.Range("G5").Value = 1667
$ws.Range("H5").Value = 1267
$ws.Range("I5").Value = 1255
$ws.Range("J5").Value = 12
$ws.Range("K5").Value = 7796
$ws.Range("L5").Value = 2631
$ws.Range("M5").Value = 5165
$ws.Range("N5").Value = 5120
$ws.Range("O5").Value = 44
$ws.Range("P5").Value = 468
$ws.Range("Q5").Value = -1183
$ws.Range("R5").Value = 708
$ws.Range("S5").Value = 714
$ws.Range("T5").Value = 4
$ws.Range("V5").Value = 1495
$ws.Range("W5").Value = 74.97
$ws.Range("X5").Value = 56.96
$ws.Range("Y5").Value = 27.12
$ws.Range("Z5").Value = 18.81
$ws.Range("AA5").Value = 50.95
$ws.Range("AB5").Value = 1026.03
$ws.Range("AC5").Value = 1013
$ws.Range("AD5").Value = 6.05
$ws.Range("AE5").Value = 4167
$ws.Range("AF5").Value = 1.61
$ws.Range("AG5").Value = 151
$ws.Range("AH5").Value = 2.46
$ws.Range("AI5").Value = 14.72
$ws.Range("AJ5").Value = 123977752
$ws.Range("U5").ClearContents()

# Row 6
$ws.Range("D6").Value = 2266
$ws.Range("E6").Value = 1379
$ws.Range("F6").Value = 1379
$ws.Range("G6").Value = 1385
$ws.Range("H6").Value = 1034
$ws.Range("I6").Value = 1024
$ws.Range("K6").Value = 10081
$ws.Range("L6").Value = 4150
$ws.Range("M6").Value = 5931
$ws.Range("N6").Value = 5931
$ws.Range("P6").Value = 514
$ws.Range("Q6").Value = -1686
$ws.Range("R6").Value = -99
$ws.Range("S6").Value = 1817
$ws.Range("T6").Value = 2
$ws.Range("V6").Value = 3460
$ws.Range("W6").Value = 60.85
$ws.Range("X6").Value = 45.63
$ws.Range("Y6").Value = 18.52
$ws.Range("Z6").Value = 11.57
$ws.Range("AA6").Value = 69.95999999999999
$ws.Range("AB6").Value = 1079.1
$ws.Range("AC6").Value = 826
$ws.Range("AD6").Value = 4.85
$ws.Range("AE6").Value = 4847
$ws.Range("AF6").Value = 0.91
$ws.Range("AG6").Value = 166
$ws.Range("AH6").Value = 4.14
$ws.Range("AI6").Value = 19.76
$ws.Range("AJ6").Value = 123977753
$ws.Range("U6").ClearContents()

# Row 7
$ws.Range("D7").Value = 2238
$ws.Range("E7").Value = 1104
$ws.Range("G7").Value = 1138
$ws.Range("H7").Value = 858
$ws.Range("I7").Value = 854
$ws.Range("K7").Value = 12896
$ws.Range("L7").Value = 6152
$ws.Range("M7").Value = 6744
$ws.Range("N7").Value = 6697
$ws.Range("P7").Value = 537
$ws.Range("W7").Value = 49.31
$ws.Range("X7").Value = 38.31
$ws.Range("Y7").Value = 13.53
$ws.Range("Z7").Value = 7.46
$ws.Range("AA7").Value = 91.23999999999999
$ws.Range("AC7").Value = 689
$ws.Range("AD7").Value = 4.27
$ws.Range("AE7").Value = 5473
$ws.Range("AF7").Value = 0.54
$ws.Range("AG7").Value = 182
$ws.Range("AH7").Value = 6.19
$ws.Range("AI7").Value = 24.04
$ws.Range("Q7").ClearContents()
$ws.Range("R7").ClearContents()
$ws.Range("S7").ClearContents()
$ws.Range("T7").ClearContents()
$ws.Range("U7").ClearContents()

# Row 8
$ws.Range("D8").Value = 2136
$ws.Range("E8").Value = 1102
$ws.Range("G8").Value = 1111
$ws.Range("H8").Value = 835
$ws.Range("I8").Value = 833
$ws.Range("K8").Value = 12860
$ws.Range("L8").Value = 5484
$ws.Range("M8").Value = 7376
$ws.Range("N8").Value = 7391
$ws.Range("P8").Value = 537
$ws.Range("W8").Value = 51.58
$ws.Range("X8").Value = 39.08
$ws.Range("Y8").Value = 11.83
$ws.Range("Z8").Value = 6.48
$ws.Range("AA8").Value = 74.34999999999999
$ws.Range("AC8").Value = 672
$ws.Range("AD8").Value = 4.38
$ws.Range("AE8").Value = 6040
$ws.Range("AF8").Value = 0.49
$ws.Range("AG8").Value = 173
$ws.Range("AH8").Value = 5.88
$ws.Range("AI8").Value = 23.43
$ws.Range("Q8").ClearContents()
$ws.Range("R8").ClearContents()
$ws.Range("S8").ClearContents()
$ws.Range("T8").ClearContents()
$ws.Range("U8").ClearContents()

# Row 9
$ws.Range("D9").Value = 2104
$ws.Range("E9").Value = 1138
$ws.Range("G9").Value = 1148
$ws.Range("H9").Value = 863
$ws.Range("I9").Value = 862
$ws.Range("K9").Value = 13377
$ws.Range("L9").Value = 5340
$ws.Range("M9").Value = 8037
$ws.Range("N9").Value = 8162
$ws.Range("P9").Value = 537
$ws.Range("W9").Value = 54.1
$ws.Range("X9").Value = 41.01
$ws.Range("Y9").Value = 11.08
$ws.Range("Z9").Value = 6.58
$ws.Range("AA9").Value = 66.44
$ws.Range("AC9").Value = 695
$ws.Range("AD9").Value = 4.23
$ws.Range("AE9").Value = 6670
$ws.Range("AF9").Value = 0.44
$ws.Range("AG9").Value = 173
$ws.Range("AH9").Value = 5.88
$ws.Range("AI9").Value = 22.65
$ws.Range("Q9").ClearContents()
$ws.Range("R9").ClearContents()
$ws.Range("S9").ClearContents()
$ws.Range("T9").ClearContents()
$ws.Range("U9").ClearContents()
